# Update the "Förändrad" (Changed) date column (C) for rows 2-6
# from 45175 (2023-09-06) to 45183 (2023-09-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
